$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "Day 2" section (rows 9-11): fill in the summary counts ---
$ws.Range("C9").Value = 456
$ws.Range("C10").Value = 703
$ws.Range("C11").Value = 511

# --- "Day 3" section (rows 17-19): fill in the summary counts ---
$ws.Range("C17").Value = 538
$ws.Range("C18").Value = 733
$ws.Range("C19").Value = 511

# --- Update the sheet view: scroll down and select C18 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
